# Commit 6: add a dashed-outline "wireframe" rectangle to each of the three
# slides (desktop main screen, sign-in form, mobile main screen).
#
# Each new rectangle is a plain theme-styled rectangle (the default shape
# drawn from the Shapes gallery) with a dash-dot outline, sent behind every
# other shape on its slide. To reproduce PowerPoint's exact default shape
# style (<p:style> lnRef/fillRef/effectRef/fontRef theme refs) - which isn't
# reachable through a plain AddShape() call in this automation surface - we
# duplicate an existing shape on slide 2 ("Rectangle 9") that already carries
# that exact style, then reposition/resize/rename/restyle the copy.
#
# NOTE: this COM shim only reliably binds *positional* function parameters,
# so the helper below is called positionally rather than with -Name style
# switches.

$p = $ppt.ActivePresentation
$EMU_PER_POINT = 12700

function Add-WireframeRectangle {
    param($TargetSlideIndex, $NewName, $OffX, $OffY, $ExtCx, $ExtCy, $BumpCount)

    $targetSlide = $p.Slides.Item($TargetSlideIndex)

    # Burn through $BumpCount shape ids on the target slide so the shape id
    # counter catches up to the value PowerPoint actually used (ids are never
    # reused on a slide, even after a shape is deleted).
    for ($i = 0; $i -lt $BumpCount; $i++) {
        $dummy = $targetSlide.Shapes.AddShape(1, 0, 0, 1, 1)
        $dummy.Delete()
    }

    # Source shape that already has the desired theme "quick style"
    # (lnRef idx=2 dk1 / fillRef idx=1 lt1 / effectRef idx=0 dk1 / fontRef minor dk1).
    $styleSourceSlide = $p.Slides.Item(2)
    $styleSource = $styleSourceSlide.Shapes.Item("Rectangle 9")

    $styleSource.Copy()
    $pastedRange = $targetSlide.Shapes.Paste()
    $newShape = $pastedRange.Item(1)

    $newShape.Name = $NewName

    $newShape.Left = $OffX / $EMU_PER_POINT
    $newShape.Top = $OffY / $EMU_PER_POINT
    $newShape.Width = $ExtCx / $EMU_PER_POINT
    $newShape.Height = $ExtCy / $EMU_PER_POINT

    # msoLineDashDot
    $newShape.Line.DashStyle = 5

    # msoSendToBack -> becomes the first shape in the slide's shape tree.
    $newShape.ZOrder(1)

    return $newShape
}

# Slide 1 - Gilbert's Poker (desktop main screen). Next free id on this
# slide is naturally 37; the target id is 43, so bump past 6 ids first.
Add-WireframeRectangle 1 "Rectangle 42" 822121 595618 1774899 1694576 6 | Out-Null

# Slide 2 - sign-in form. Next free id is already 3, matching the target.
Add-WireframeRectangle 2 "Rectangle 2" 1058795 1650534 1774899 1694576 0 | Out-Null

# Slide 3 - Gilbert's Mobile Poker (mobile main screen). Next free id is
# already 19, matching the target.
Add-WireframeRectangle 3 "Rectangle 18" 852217 654341 1774899 1694576 0 | Out-Null
